# Apply "toy-spam min 5" edit: update word-frequency tables and
# remove rows/cells that no longer meet the minimum-count-5 threshold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (negative word list shrank from 32 to 31 entries)
$ws.Range("A34:H34").EntireRow.Delete()

# Remove the last entry of the positive word list (shrank from 19 to 18 entries)
$ws.Range("J21:Q21").Clear()

# Update all changed cell values to reflect the recomputed word-frequency tables
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 46
$ws.Range("D3").Value = 46
$ws.Range("H3").Value = 0
$ws.Range("K4").Value = 0.8307692307692308
$ws.Range("L4").Value = 54
$ws.Range("M4").Value = 54
$ws.Range("Q4").Value = 11
$ws.Range("A5").Value = "returned"
$ws.Range("B5").Value = 0.7631578947368421
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = 29
$ws.Range("H5").Value = 9
$ws.Range("K5").Value = 0.7956989247311828
$ws.Range("L5").Value = 74
$ws.Range("M5").Value = 74
$ws.Range("Q5").Value = 19
$ws.Range("A6").Value = "poor"
$ws.Range("B6").Value = 0.7464788732394366
$ws.Range("C6").Value = 53
$ws.Range("D6").Value = 53
$ws.Range("H6").Value = 18
$ws.Range("K6").Value = 0.78125
$ws.Range("L6").Value = 50
$ws.Range("M6").Value = 50
$ws.Range("Q6").Value = 14
$ws.Range("J7").Value = "classic"
$ws.Range("K7").Value = 0.6981132075471698
$ws.Range("L7").Value = 37
$ws.Range("M7").Value = 37
$ws.Range("Q7").Value = 16
$ws.Range("A8").Value = "disappointed"
$ws.Range("B8").Value = 0.6666666666666666
$ws.Range("C8").Value = 124
$ws.Range("D8").Value = 124
$ws.Range("H8").Value = 62
$ws.Range("J8").Value = "thank"
$ws.Range("K8").Value = 0.6231884057971014
$ws.Range("L8").Value = 43
$ws.Range("M8").Value = 43
$ws.Range("Q8").Value = 26
$ws.Range("A9").Value = "waste"
$ws.Range("B9").Value = 0.6216216216216216
$ws.Range("C9").Value = 92
$ws.Range("D9").Value = 92
$ws.Range("H9").Value = 56
$ws.Range("K9").Value = 0.5437589670014347
$ws.Range("L9").Value = 379
$ws.Range("M9").Value = 379
$ws.Range("Q9").Value = 318
$ws.Range("A10").Value = "junk"
$ws.Range("B10").Value = 0.6181818181818182
$ws.Range("C10").Value = 34
$ws.Range("D10").Value = 34
$ws.Range("H10").Value = 21
$ws.Range("K10").Value = 0.479253112033195
$ws.Range("L10").Value = 231
$ws.Range("M10").Value = 231
$ws.Range("Q10").Value = 251
$ws.Range("A11").Value = "broke"
$ws.Range("B11").Value = 0.6067961165048543
$ws.Range("C11").Value = 125
$ws.Range("D11").Value = 125
$ws.Range("H11").Value = 81
$ws.Range("K11").Value = 0.4491803278688525
$ws.Range("L11").Value = 548
$ws.Range("M11").Value = 548
$ws.Range("Q11").Value = 672
$ws.Range("A12").Value = "smaller"
$ws.Range("B12").Value = 0.5966386554621849
$ws.Range("C12").Value = 71
$ws.Range("D12").Value = 71
$ws.Range("H12").Value = 48
$ws.Range("K12").Value = 0.3577981651376147
$ws.Range("L12").Value = 117
$ws.Range("M12").Value = 117
$ws.Range("Q12").Value = 210
$ws.Range("A13").Value = "small"
$ws.Range("B13").Value = 0.5072463768115942
$ws.Range("C13").Value = 175
$ws.Range("D13").Value = 175
$ws.Range("H13").Value = 170
$ws.Range("K13").Value = 0.3333333333333333
$ws.Range("L13").Value = 40
$ws.Range("M13").Value = 40
$ws.Range("Q13").Value = 80
$ws.Range("B14").Value = 0.4216867469879518
$ws.Range("C14").Value = 35
$ws.Range("D14").Value = 35
$ws.Range("H14").Value = 48
$ws.Range("K14").Value = 0.3192771084337349
$ws.Range("L14").Value = 53
$ws.Range("M14").Value = 53
$ws.Range("Q14").Value = 113
$ws.Range("A15").Value = "plastic"
$ws.Range("B15").Value = 0.4015748031496063
$ws.Range("C15").Value = 51
$ws.Range("D15").Value = 51
$ws.Range("H15").Value = 76
$ws.Range("K15").Value = 0.3121693121693122
$ws.Range("L15").Value = 59
$ws.Range("M15").Value = 59
$ws.Range("Q15").Value = 130
$ws.Range("A16").Value = "ok"
$ws.Range("B16").Value = 0.390625
$ws.Range("C16").Value = 50
$ws.Range("D16").Value = 50
$ws.Range("H16").Value = 78
$ws.Range("J16").Value = "happy"
$ws.Range("K16").Value = 0.2027972027972028
$ws.Range("L16").Value = 29
$ws.Range("M16").Value = 29
$ws.Range("Q16").Value = 114
$ws.Range("A17").Value = "cheap"
$ws.Range("B17").Value = 0.3696682464454976
$ws.Range("C17").Value = 78
$ws.Range("D17").Value = 78
$ws.Range("H17").Value = 133
$ws.Range("K17").Value = 0.1935483870967742
$ws.Range("L17").Value = 36
$ws.Range("M17").Value = 36
$ws.Range("Q17").Value = 150
$ws.Range("A18").Value = "apart"
$ws.Range("B18").Value = 0.3473684210526316
$ws.Range("C18").Value = 33
$ws.Range("D18").Value = 33
$ws.Range("H18").Value = 62
$ws.Range("J18").Value = "christmas"
$ws.Range("K18").Value = 0.1686746987951807
$ws.Range("L18").Value = 42
$ws.Range("M18").Value = 42
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 207
$ws.Range("A19").Value = "difficult"
$ws.Range("B19").Value = 0.3258426966292135
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 29
$ws.Range("H19").Value = 60
$ws.Range("J19").Value = "fun"
$ws.Range("K19").Value = 0.1578947368421053
$ws.Range("L19").Value = 180
$ws.Range("M19").Value = 181
$ws.Range("N19").Value = 0.99
$ws.Range("O19").Value = 0.01000000000000001
$ws.Range("P19").Value = $true
$ws.Range("Q19").Value = 960
$ws.Range("A20").Value = "thought"
$ws.Range("B20").Value = 0.301980198019802
$ws.Range("C20").Value = 61
$ws.Range("D20").Value = 61
$ws.Range("H20").Value = 141
$ws.Range("K20").Value = 0.0792722547108512
$ws.Range("L20").Value = 122
$ws.Range("M20").Value = 124
$ws.Range("N20").Value = 0.98
$ws.Range("O20").Value = 0.02000000000000002
$ws.Range("Q20").Value = 1417
$ws.Range("A21").Value = "size"
$ws.Range("B21").Value = 0.2010309278350516
$ws.Range("H21").Value = 155
$ws.Range("A22").Value = "hard"
$ws.Range("B22").Value = 0.175
$ws.Range("C22").Value = 35
$ws.Range("D22").Value = 35
$ws.Range("H22").Value = 165
$ws.Range("A23").Value = "work"
$ws.Range("B23").Value = 0.1746031746031746
$ws.Range("C23").Value = 55
$ws.Range("D23").Value = 56
$ws.Range("E23").Value = 0.02
$ws.Range("F23").Value = 0.98
$ws.Range("G23").Value = $true
$ws.Range("H23").Value = 260
$ws.Range("A24").Value = "would"
$ws.Range("B24").Value = 0.1711309523809524
$ws.Range("C24").Value = 115
$ws.Range("D24").Value = 117
$ws.Range("E24").Value = 0.02
$ws.Range("F24").Value = 0.98
$ws.Range("G24").Value = $true
$ws.Range("H24").Value = 557
$ws.Range("A25").Value = "item"
$ws.Range("B25").Value = 0.1666666666666667
$ws.Range("C25").Value = 46
$ws.Range("D25").Value = 46
$ws.Range("H25").Value = 230
$ws.Range("A26").Value = "money"
$ws.Range("B26").Value = 0.1645569620253164
$ws.Range("C26").Value = 52
$ws.Range("D26").Value = 52
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 264
$ws.Range("A27").Value = "better"
$ws.Range("B27").Value = 0.1448598130841121
$ws.Range("C27").Value = 31
$ws.Range("D27").Value = 31
$ws.Range("H27").Value = 183
$ws.Range("A28").Value = "price"
$ws.Range("B28").Value = 0.1440922190201729
$ws.Range("C28").Value = 50
$ws.Range("D28").Value = 51
$ws.Range("E28").Value = 0.02
$ws.Range("F28").Value = 0.98
$ws.Range("G28").Value = $true
$ws.Range("H28").Value = 297
$ws.Range("A29").Value = "product"
$ws.Range("B29").Value = 0.1280353200883002
$ws.Range("C29").Value = 58
$ws.Range("D29").Value = 59
$ws.Range("E29").Value = 0.02
$ws.Range("F29").Value = 0.98
$ws.Range("G29").Value = $true
$ws.Range("H29").Value = 395
$ws.Range("A30").Value = "little"
$ws.Range("B30").Value = 0.08258928571428571
$ws.Range("C30").Value = 37
$ws.Range("D30").Value = 38
$ws.Range("E30").Value = 0.03
$ws.Range("F30").Value = 0.97
$ws.Range("G30").Value = $true
$ws.Range("H30").Value = 411
$ws.Range("A31").Value = "use"
$ws.Range("B31").Value = 0.07967032967032966
$ws.Range("C31").Value = 29
$ws.Range("D31").Value = 30
$ws.Range("H31").Value = 335
$ws.Range("B32").Value = 0.06765676567656766
$ws.Range("C32").Value = 41
$ws.Range("D32").Value = 43
$ws.Range("H32").Value = 565
$ws.Range("A33").Value = "one"
$ws.Range("B33").Value = 0.0457433290978399
$ws.Range("C33").Value = 36
$ws.Range("D33").Value = 43
$ws.Range("E33").Value = 0.16
$ws.Range("F33").Value = 0.84
$ws.Range("H33").Value = 751
